$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, even if it looks like a number,
# without permanently altering the cell style (e.g. "303.70" must stay text).
function Set-TextValue($ws, $ref, $text) {
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = '43.173.15'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue $ws 'D5' '303.70'
Set-TextValue $ws 'D6' '97.88'
$ws.Range("E6").Value = '  +0.75%  '
Set-TextValue $ws 'D7' '0.507'
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.37%  '
Set-TextValue $ws 'D10' '35.65'
$ws.Range("E10").Value = '  +0.16%  '
Set-TextValue $ws 'D11' '19.33'
$ws.Range("E11").Value = '  +8.01%  '
Set-TextValue $ws 'D12' '0.0792'
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("E14").Value = '  +2.05%  '
$ws.Range("D15").Value = '2.689.79'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").Value = '2.317.24'
$ws.Range("E16").Value = '  +0.76%  '
Set-TextValue $ws 'D17' '0.786'
$ws.Range("E17").Value = '  +1.21%  '
$ws.Range("D18").Value = '43.092.94'
$ws.Range("E18").Value = '  +0.59%  '
Set-TextValue $ws 'D19' '12.59'
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("E21").Value = '  +0.83%  '
Set-TextValue $ws 'D22' '67.98'
$ws.Range("E22").Value = '  +0.20%  '
Set-TextValue $ws 'D23' '237.95'
$ws.Range("E23").Value = '  -0.94%  '
Set-TextValue $ws 'D24' '2.21'
$ws.Range("E24").Value = '  +3.84%  '
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  +0.20%  '
Set-TextValue $ws 'D27' '24.99'
$ws.Range("E27").Value = '  -1.37%  '
Set-TextValue $ws 'D28' '2.39'
$ws.Range("E28").Value = '  +18.26%  '
Set-TextValue $ws 'D29' '165.99'
$ws.Range("E29").Value = '  +0.45%  '
Set-TextValue $ws 'D30' '9.12'
$ws.Range("E30").Value = '  +0.91%  '
Set-TextValue $ws 'D31' '33.11'
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("E32").Value = '  -0.03%  '
Set-TextValue $ws 'D33' '18.06'
$ws.Range("E33").Value = '  +6.51%  '
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("E35").Value = '  -8.32%  '
Set-TextValue $ws 'D37' '0.0694'
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("E38").Value = '  +0.38%  '
Set-TextValue $ws 'D39' '2.79'
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").Value = '2.000.29'
$ws.Range("E42").Value = '  -0.60%  '
Set-TextValue $ws 'D43' '10.70'
$ws.Range("E43").Value = '  +5.31%  '
$ws.Range("E44").Value = '  +0.27%  '
Set-TextValue $ws 'D45' '18.21'
$ws.Range("E45").Value = '  +4.56%  '
$ws.Range("E46").Value = '  -1.57%  '
Set-TextValue $ws 'D47' '2.78'
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D48").Value = '2.556.60'
$ws.Range("E48").Value = '  +1.20%  '
Set-TextValue $ws 'D49' '53.71'
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("E50").Value = '  -5.87%  '
Set-TextValue $ws 'D51' '72.07'
$ws.Range("E51").Value = '  -0.03%  '
